# Auto-generated cell value updates (Sheets market-price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 86
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -877
$ws.Range("N86").ClearContents()
# row 89
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4384
$ws.Range("N89").ClearContents()
# row 100
$ws.Range("H100").Value = 40001860
$ws.Range("I100").Value = 40001860
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 40001860
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -40001319
$ws.Range("N100").ClearContents()
# row 112
$ws.Range("H112").Value = 1312.9464
$ws.Range("J112").Value = 1312.9464
$ws.Range("L112").Value = 3938.8392
$ws.Range("N112").Value = -6154.8392
# row 129
$ws.Range("H129").Value = 837.7041
$ws.Range("I129").Value = 319.4
$ws.Range("J129").Value = 865.5699
$ws.Range("K129").Value = 958.1999999999999
$ws.Range("L129").Value = 2596.7097
$ws.Range("M129").Value = 4041.8
$ws.Range("N129").Value = -12596.7097
# row 135
$ws.Range("H135").Value = 1068.3334
$ws.Range("I135").Value = 728.8
$ws.Range("J135").Value = 1492.75
$ws.Range("K135").Value = 6559.2
$ws.Range("L135").Value = 13434.75
$ws.Range("M135").Value = -4024.2
$ws.Range("N135").Value = -18504.75

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 692.625
$ws.Range("I2").Value = 692.625
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 692.625
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -579.625
$ws.Range("N2").ClearContents()
# row 23
$ws.Range("H23").Value = 26388.75
$ws.Range("J23").Value = 26388.75
$ws.Range("L23").Value = 26388.75
$ws.Range("N23").Value = -26906.75
# row 31
$ws.Range("H31").Value = 22500
$ws.Range("I31").Value = 5000
$ws.Range("J31").Value = 40000
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 40000
$ws.Range("M31").Value = -4706
$ws.Range("N31").Value = -40588
# row 32
$ws.Range("H32").Value = 7779.423
$ws.Range("I32").Value = 7929.75
$ws.Range("J32").Value = 7650.5713
$ws.Range("K32").Value = 7929.75
$ws.Range("L32").Value = 7650.5713
$ws.Range("M32").Value = -7642.75
$ws.Range("N32").Value = -8224.5713
# row 61
$ws.Range("H61").Value = 1734.9333
$ws.Range("I61").Value = 1461.4615
$ws.Range("K61").Value = 1461.4615
$ws.Range("M61").Value = -1249.4615
# row 74
$ws.Range("H74").Value = 4564.231
$ws.Range("I74").Value = 4203.8887
$ws.Range("J74").Value = 5375
$ws.Range("K74").Value = 4203.8887
$ws.Range("L74").Value = 5375
$ws.Range("M74").Value = -3329.8887
$ws.Range("N74").Value = -7123
# row 77
$ws.Range("H77").Value = 4564.231
$ws.Range("I77").Value = 4203.8887
$ws.Range("J77").Value = 5375
$ws.Range("K77").Value = 21019.4435
$ws.Range("L77").Value = 26875
$ws.Range("M77").Value = -16651.4435
$ws.Range("N77").Value = -35611
# row 97
$ws.Range("H97").Value = 1320.1666
$ws.Range("I97").Value = 1184.2
$ws.Range("K97").Value = 1184.2
$ws.Range("M97").Value = -688.2
# row 110
$ws.Range("H110").Value = 666.3
$ws.Range("I110").Value = 596.4286
$ws.Range("K110").Value = 596.4286
$ws.Range("M110").Value = 1448.5714
# row 116
$ws.Range("H116").Value = 692.625
$ws.Range("I116").Value = 692.625
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 692.625
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1601.375
$ws.Range("N116").ClearContents()
# row 136
$ws.Range("H136").Value = 1734.9333
$ws.Range("I136").Value = 1461.4615
$ws.Range("K136").Value = 4384.3845
$ws.Range("M136").Value = -1834.3845
# row 139
$ws.Range("H139").Value = 41587.414
$ws.Range("J139").Value = 41587.414
$ws.Range("L139").Value = 41587.414
$ws.Range("N139").Value = -51867.414

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 692.625
$ws.Range("I3").Value = 692.625
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 692.625
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -578.625
$ws.Range("N3").ClearContents()
# row 26
$ws.Range("H26").Value = 25000
$ws.Range("I26").Value = 25000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -24708
$ws.Range("N26").ClearContents()
# row 99
$ws.Range("H99").Value = 3287.8333
$ws.Range("I99").Value = 1400.1111
$ws.Range("K99").Value = 1400.1111
$ws.Range("M99").Value = 97.88889999999992
# row 134
$ws.Range("H134").Value = 2546.6287
$ws.Range("I134").Value = 1288.9286
$ws.Range("K134").Value = 3866.7858
$ws.Range("M134").Value = -1331.7858

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 5219.7393
$ws.Range("I31").Value = 2141.0322
$ws.Range("J31").Value = 11582.4
$ws.Range("K31").Value = 2141.0322
$ws.Range("L31").Value = 11582.4
$ws.Range("M31").Value = -1846.0322
$ws.Range("N31").Value = -12172.4
# row 34
$ws.Range("H34").Value = 5219.7393
$ws.Range("I34").Value = 2141.0322
$ws.Range("J34").Value = 11582.4
$ws.Range("K34").Value = 2141.0322
$ws.Range("L34").Value = 11582.4
$ws.Range("M34").Value = -1939.0322
$ws.Range("N34").Value = -11986.4
# row 138
$ws.Range("H138").Value = 41328.332
$ws.Range("J138").Value = 41328.332
$ws.Range("L138").Value = 41328.332
$ws.Range("N138").Value = -51608.332
# row 140
$ws.Range("H140").Value = 105326.664
$ws.Range("J140").Value = 117242.5
$ws.Range("L140").Value = 117242.5
$ws.Range("N140").Value = -127602.5
# row 141
$ws.Range("H141").Value = 30931.25
$ws.Range("J141").Value = 30931.25
$ws.Range("L141").Value = 30931.25
$ws.Range("N141").Value = -41291.25

$ws = $wb.Worksheets.Item("CUL")
# row 59
$ws.Range("H59").Value = 2899.5

$ws = $wb.Worksheets.Item("GSM")
# row 113
$ws.Range("H113").Value = 2310.9375
$ws.Range("I113").Value = 2851.375
$ws.Range("J113").Value = 1770.5
$ws.Range("K113").Value = 2851.375
$ws.Range("L113").Value = 1770.5
$ws.Range("M113").Value = -681.375
$ws.Range("N113").Value = -6110.5
# row 140
$ws.Range("H140").Value = 38503.527
$ws.Range("J140").Value = 38503.527
$ws.Range("L140").Value = 38503.527
$ws.Range("N140").Value = -48863.527
# row 141
$ws.Range("H141").Value = 41167.5
$ws.Range("J141").Value = 42890
$ws.Range("L141").Value = 42890
$ws.Range("N141").Value = -53250

$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 7801.909
$ws.Range("I132").Value = 3756.8
$ws.Range("J132").Value = 8991.647000000001
$ws.Range("K132").Value = 11270.4
$ws.Range("L132").Value = 26974.941
$ws.Range("M132").Value = -8740.400000000001
$ws.Range("N132").Value = -32034.941
# row 136
$ws.Range("H136").Value = 5200.5
$ws.Range("I136").Value = 1634.8334
$ws.Range("J136").Value = 6983.3335
$ws.Range("K136").Value = 4904.5002
$ws.Range("L136").Value = 20950.0005
$ws.Range("M136").Value = -2354.5002
$ws.Range("N136").Value = -26050.0005
# row 140
$ws.Range("H140").Value = 66214.08
$ws.Range("J140").Value = 66214.08
$ws.Range("L140").Value = 66214.08
$ws.Range("N140").Value = -76574.08
# row 141
$ws.Range("H141").Value = 32125.79
$ws.Range("J141").Value = 32125.79
$ws.Range("L141").Value = 32125.79
$ws.Range("N141").Value = -42485.79

$ws = $wb.Worksheets.Item("WVR")
# row 46
$ws.Range("H46").Value = 46616.08
$ws.Range("J46").Value = 46616.08
$ws.Range("L46").Value = 46616.08
$ws.Range("N46").Value = -47078.08
# row 122
$ws.Range("H122").Value = 6107.96
$ws.Range("I122").Value = 4208.6875
$ws.Range("J122").Value = 9484.444
$ws.Range("K122").Value = 12626.0625
$ws.Range("L122").Value = 28453.332
$ws.Range("M122").Value = -10176.0625
$ws.Range("N122").Value = -33353.33199999999
# row 132
$ws.Range("H132").Value = 27789270
$ws.Range("I132").Value = 26175
$ws.Range("K132").Value = 78525
$ws.Range("M132").Value = -75995
# row 134
$ws.Range("H134").Value = 46616.08
$ws.Range("J134").Value = 46616.08
$ws.Range("L134").Value = 139848.24
$ws.Range("N134").Value = -144918.24
# row 135
$ws.Range("H135").Value = 41715
$ws.Range("J135").Value = 41715
$ws.Range("L135").Value = 41715
$ws.Range("N135").Value = -51855
# row 136
$ws.Range("H136").Value = 8930.218000000001
$ws.Range("I136").Value = 10060.833
$ws.Range("J136").Value = 7696.8184
$ws.Range("K136").Value = 30182.499
$ws.Range("L136").Value = 23090.4552
$ws.Range("M136").Value = -27632.499
$ws.Range("N136").Value = -28190.4552
# row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
# row 140
$ws.Range("H140").Value = 32275.334
$ws.Range("J140").Value = 32275.334
$ws.Range("L140").Value = 32275.334
$ws.Range("N140").Value = -42635.334
